$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "БИВТ-22-17": raise the exam-pass threshold from 35 to 40, and fill
# in a few lab scores that were missing / recorded as text.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("БИВТ-22-17")

$ws1.Range("O2").Formula  = "=IF(M2>=40,""pass"",""fail"")"
$ws1.Range("O3").Formula  = "=IF(M3>=40,""pass"",""fail"")"
$ws1.Range("O4:O31").Formula = "=IF(M4>=40,""pass"",""fail"")"
$ws1.Range("O32").Formula = "=IF(M32>=40,""pass"",""fail"")"

$ws1.Range("I23").Value = 5
$ws1.Range("E25").Value = 5
$ws1.Range("F25").Value = 5

[void]$ws1.Activate()
[void]$ws1.Range("O32").Select()

# ---------------------------------------------------------------------------
# Sheet "БИВТ-22-18": same threshold change, plus several lab scores
# recorded as "pass" text were corrected to the numeric score 5, and a
# couple of blank scores were filled in.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("БИВТ-22-18")

$ws2.Range("O2").Formula  = "=IF(M2>=40,""pass"",""fail"")"
$ws2.Range("O3").Formula  = "=IF(M3>=40,""pass"",""fail"")"
$ws2.Range("O4:O25").Formula = "=IF(M4>=40,""pass"",""fail"")"
$ws2.Range("O26").Formula = "=IF(M26>=40,""pass"",""fail"")"

$ws2.Range("D4").Value = 5
$ws2.Range("H7").Value = 5
$ws2.Range("B10").Value = 5
$ws2.Range("C10").Value = 2.9
$ws2.Range("D11").Value = 5
$ws2.Range("E11").Value = 5
$ws2.Range("D19").Value = 5
$ws2.Range("E20").Value = 5
$ws2.Range("D22").Value = 5
$ws2.Range("F23").Value = 5

[void]$ws2.Activate()
[void]$ws2.Range("O25").Select()

# ---------------------------------------------------------------------------
# Sheet "БИВТ-22-20": same threshold change only.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("БИВТ-22-20")

$ws3.Range("O2").Formula  = "=IF(M2>=40,""pass"",""fail"")"
$ws3.Range("O3").Formula  = "=IF(M3>=40,""pass"",""fail"")"
$ws3.Range("O4:O29").Formula = "=IF(M4>=40,""pass"",""fail"")"
$ws3.Range("O30").Formula = "=IF(M30>=40,""pass"",""fail"")"

[void]$ws3.Activate()
[void]$ws3.Range("O31").Select()

# Leave the workbook on the first sheet, as in the original file.
[void]$ws1.Activate()
